$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value2 = 2933.5
$ws.Range("I69").Value2 = 2200
$ws.Range("K69").Value2 = 6600
$ws.Range("M69").Value2 = -5726
# Row 72
$ws.Range("H72").Value2 = 2933.5
$ws.Range("I72").Value2 = 2200
$ws.Range("K72").Value2 = 19800
$ws.Range("M72").Value2 = -15432
# Row 86
$ws.Range("H86").Value2 = 1686.6666
$ws.Range("I86").Value2 = 1560.4
$ws.Range("J86").Value2 = 1749.8
$ws.Range("K86").Value2 = 1560.4
$ws.Range("L86").Value2 = 1749.8
$ws.Range("M86").Value2 = -437.4000000000001
$ws.Range("N86").Value2 = -3995.8
# Row 89
$ws.Range("H89").Value2 = 1686.6666
$ws.Range("I89").Value2 = 1560.4
$ws.Range("J89").Value2 = 1749.8
$ws.Range("K89").Value2 = 7802
$ws.Range("L89").Value2 = 8749
$ws.Range("M89").Value2 = -2186
$ws.Range("N89").Value2 = -19981
# Row 129
$ws.Range("H129").Value2 = 659.44446
$ws.Range("I129").Value2 = 287
$ws.Range("J129").Value2 = 1125
$ws.Range("K129").Value2 = 861
$ws.Range("L129").Value2 = 3375
$ws.Range("M129").Value2 = 4139
$ws.Range("N129").Value2 = -13375
# Row 132
$ws.Range("H132").Value2 = 629574.6
$ws.Range("I132").Value2 = 1266.6
$ws.Range("J132").Value2 = 3771115
$ws.Range("K132").Value2 = 3799.8
$ws.Range("L132").Value2 = 11313345
$ws.Range("M132").Value2 = -1269.8
$ws.Range("N132").Value2 = -11318405
# Row 138
$ws.Range("H138").Value2 = 1685344.1
$ws.Range("I138").Value2 = 1085.0613
$ws.Range("K138").Value2 = 3255.1839
$ws.Range("M138").Value2 = 1884.8161

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value2 = 23200600
$ws.Range("I8").Value2 = 29000000
$ws.Range("J8").Value2 = 3000
$ws.Range("K8").Value2 = 29000000
$ws.Range("L8").Value2 = 3000
$ws.Range("M8").Value2 = -28999856
$ws.Range("N8").Value2 = -3288
# Row 32
$ws.Range("H32").Value2 = 772.83
$ws.Range("I32").Value2 = 680.69476
$ws.Range("J32").Value2 = 2523.4
$ws.Range("K32").Value2 = 680.69476
$ws.Range("L32").Value2 = 2523.4
$ws.Range("M32").Value2 = -393.69476
$ws.Range("N32").Value2 = -3097.4
# Row 44
$ws.Range("H44").Value2 = 15000
# Row 74
$ws.Range("H74").Value2 = 6147429
$ws.Range("I74").Value2 = 7606885.5
$ws.Range("J74").Value2 = 127168.75
$ws.Range("K74").Value2 = 7606885.5
$ws.Range("L74").Value2 = 127168.75
$ws.Range("M74").Value2 = -7606011.5
$ws.Range("N74").Value2 = -128916.75
# Row 77
$ws.Range("H77").Value2 = 6147429
$ws.Range("I77").Value2 = 7606885.5
$ws.Range("J77").Value2 = 127168.75
$ws.Range("K77").Value2 = 38034427.5
$ws.Range("L77").Value2 = 635843.75
$ws.Range("M77").Value2 = -38030059.5
$ws.Range("N77").Value2 = -644579.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value2 = 10000
$ws.Range("J92").Value2 = 10000
$ws.Range("L92").Value2 = 10000
$ws.Range("N92").Value2 = -14992
# Row 99
$ws.Range("H99").Value2 = 952.38464
$ws.Range("I99").Value2 = 818.8889
$ws.Range("J99").Value2 = 1252.75
$ws.Range("K99").Value2 = 818.8889
$ws.Range("L99").Value2 = 1252.75
$ws.Range("M99").Value2 = 679.1111
$ws.Range("N99").Value2 = -4248.75
# Row 134
$ws.Range("H134").Value2 = 1704.4255
$ws.Range("I134").Value2 = 1061.2646
$ws.Range("J134").Value2 = 3386.5386
$ws.Range("K134").Value2 = 3183.7938
$ws.Range("L134").Value2 = 10159.6158
$ws.Range("M134").Value2 = -648.7937999999999
$ws.Range("N134").Value2 = -15229.6158

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 2273.4822
$ws.Range("I31").Value2 = 1122.5135
$ws.Range("J31").Value2 = 4514.8423
$ws.Range("K31").Value2 = 1122.5135
$ws.Range("L31").Value2 = 4514.8423
$ws.Range("M31").Value2 = -827.5135
$ws.Range("N31").Value2 = -5104.8423
# Row 34
$ws.Range("H34").Value2 = 2273.4822
$ws.Range("I34").Value2 = 1122.5135
$ws.Range("J34").Value2 = 4514.8423
$ws.Range("K34").Value2 = 1122.5135
$ws.Range("L34").Value2 = 4514.8423
$ws.Range("M34").Value2 = -920.5135
$ws.Range("N34").Value2 = -4918.8423
# Row 58
$ws.Range("H58").Value2 = 14707020
$ws.Range("I58").Value2 = 16130049
$ws.Range("K58").Value2 = 16130049
$ws.Range("M58").Value2 = -16129846
# Row 62
$ws.Range("H62").Value2 = 2886.75
$ws.Range("I62").Value2 = 2417.6
$ws.Range("J62").Value2 = 3668.6667
$ws.Range("K62").Value2 = 2417.6
$ws.Range("L62").Value2 = 3668.6667
$ws.Range("M62").Value2 = -1793.6
$ws.Range("N62").Value2 = -4916.6667
# Row 65
$ws.Range("H65").Value2 = 2886.75
$ws.Range("I65").Value2 = 2417.6
$ws.Range("J65").Value2 = 3668.6667
$ws.Range("K65").Value2 = 12088
$ws.Range("L65").Value2 = 18343.3335
$ws.Range("M65").Value2 = -8968
$ws.Range("N65").Value2 = -24583.3335
# Row 136
$ws.Range("H136").Value2 = 14707020
$ws.Range("I136").Value2 = 16130049
$ws.Range("K136").Value2 = 48390147
$ws.Range("M136").Value2 = -48387597

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value2 = 1094.2979
$ws.Range("I131").Value2 = 526.55554
$ws.Range("J131").Value2 = 1228.7632
$ws.Range("K131").Value2 = 1579.66662
$ws.Range("L131").Value2 = 3686.2896
$ws.Range("M131").Value2 = 3460.33338
$ws.Range("N131").Value2 = -13766.2896
# Row 138
$ws.Range("H138").Value2 = 3208.4211
$ws.Range("I138").Value2 = 2426
$ws.Range("K138").Value2 = 7278
$ws.Range("M138").Value2 = -2138

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value2 = 500785.72
$ws.Range("I132").Value2 = 500000
$ws.Range("J132").Value2 = 502750
$ws.Range("K132").Value2 = 1500000
$ws.Range("L132").Value2 = 1508250
$ws.Range("M132").Value2 = -1497470
$ws.Range("N132").Value2 = -1513310
# Row 136
$ws.Range("H136").Value2 = 71853.80499999999
$ws.Range("I136").Value2 = 58762.156
$ws.Range("J136").Value2 = 92582.25
$ws.Range("K136").Value2 = 176286.468
$ws.Range("L136").Value2 = 277746.75
$ws.Range("M136").Value2 = -173736.468
$ws.Range("N136").Value2 = -282846.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value2 = 3451.5
$ws.Range("I122").Value2 = 2978
$ws.Range("J122").Value2 = 3640.9
$ws.Range("K122").Value2 = 8934
$ws.Range("L122").Value2 = 10922.7
$ws.Range("M122").Value2 = -6484
$ws.Range("N122").Value2 = -15822.7
# Row 132
$ws.Range("H132").Value2 = 57094.4
$ws.Range("I132").Value2 = 41500.04
$ws.Range("J132").Value2 = 184448.33
$ws.Range("K132").Value2 = 124500.12
$ws.Range("L132").Value2 = 553344.99
$ws.Range("M132").Value2 = -121970.12
$ws.Range("N132").Value2 = -558404.99
